$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/string updates (Coin name, Link, Volume%) ---
$normalUpdates = @{
    "E2" = "  +0.68%  "
    "E3" = "  +0.16%  "
    "E4" = "  +0.04%  "
    "E5" = "  +2.58%  "
    "E6" = "  +1.57%  "
    "E7" = "  +1.91%  "
    "E9" = "  +0.08%  "
    "E10" = "  -0.31%  "
    "E11" = "  +2.31%  "
    "E12" = "  +0.37%  "
    "E13" = "  +0.30%  "
    "E14" = "  -2.38%  "
    "E15" = "  +0.60%  "
    "E16" = "  -0.22%  "
    "B17" = "WrappedEther"
    "C17" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "E17" = "  +0.90%  "
    "B18" = "ShibaInu"
    "C18" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "E18" = "  -0.27%  "
    "E19" = "  -2.16%  "
    "E20" = "  -1.91%  "
    "E21" = "  -1.99%  "
    "E22" = "  -1.71%  "
    "B23" = "Dai"
    "C23" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "E23" = "  -0.10%  "
    "B24" = "Litecoin"
    "C24" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "E24" = "  -1.59%  "
    "E25" = "  +0.00%  "
    "E26" = "  +0.64%  "
    "E27" = "  -0.30%  "
    "E28" = "  +6.12%  "
    "E29" = "  +0.15%  "
    "E30" = "  +4.38%  "
    "E31" = "  +0.01%  "
    "E32" = "  -0.83%  "
    "E33" = "  -0.59%  "
    "E34" = "  +0.02%  "
    "E35" = "  +0.39%  "
    "E36" = "  -1.48%  "
    "E37" = "  -0.59%  "
    "E38" = "  +0.55%  "
    "E39" = "  -1.78%  "
    "E40" = "  +2.34%  "
    "E41" = "  -0.32%  "
    "E42" = "  -2.05%  "
    "E43" = "  -0.40%  "
    "E44" = "  -4.20%  "
    "E45" = "  -1.73%  "
    "E46" = "  -2.52%  "
    "E47" = "  +1.23%  "
    "E48" = "  -2.67%  "
    "E49" = "  -2.17%  "
    "E50" = "  -3.17%  "
    "E51" = "  -0.19%  "
}
foreach ($ref in $normalUpdates.Keys) {
    $ws.Range($ref).Value = $normalUpdates[$ref]
}

# --- Price column updates: force text storage so values like
# "0.0270" / "66.400.98" keep their exact literal formatting
# instead of being auto-coerced into floating point numbers. ---
$priceUpdates = @{
    "D2" = "66.400.98"
    "D3" = "3.297.97"
    "D5" = "588.21"
    "D6" = "180.35"
    "D9" = "3.294.38"
    "D10" = "0.126"
    "D11" = "6.86"
    "D12" = "0.402"
    "D13" = "3.878.97"
    "D15" = "66.434.72"
    "D16" = "26.56"
    "D17" = "3.302.75"
    "D18" = "0.0000163"
    "D19" = "427.12"
    "D20" = "5.48"
    "D21" = "13.00"
    "D22" = "7.28"
    "D23" = "1.00"
    "D24" = "71.38"
    "D25" = "5.69"
    "D26" = "3.457.26"
    "D27" = "0.511"
    "D28" = "0.206"
    "D29" = "0.0000114"
    "D30" = "9.33"
    "D32" = "1.92"
    "D33" = "22.28"
    "D35" = "5.17"
    "D36" = "6.56"
    "D37" = "1.18"
    "D38" = "159.12"
    "D39" = "1.43"
    "D40" = "2.850.60"
    "D41" = "1.78"
    "D42" = "26.29"
    "D43" = "4.33"
    "D44" = "0.749"
    "D45" = "39.67"
    "D46" = "5.92"
    "D47" = "2.31"
    "D48" = "0.0640"
    "D49" = "312.92"
    "D50" = "22.77"
    "D51" = "0.0270"
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$ref]
    $cell.Style = "Normal"
}
